$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data1")

# Update the tweaked input data values (Measles coverage column C)
$ws.Range("C24").Value = 0.7
$ws.Range("C25").Value = 0.68
$ws.Range("C26").Value = 0.75

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("E26").Select()
